# Updates the betting-odds worksheet:
#  - Swaps/rotates the match rows listed below back into the correct
#    chronological order (the "home"/"away" match data in columns F:V got
#    shuffled relative to the index/date columns A:E, so we restore the
#    correct pairing by swapping or rotating the F:V payload between rows).
#  - Appends one new match (Recanatese vs Pineto) as row 147.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowData {
    param($rowA, $rowB)
    $rngA = $ws.Range("F$rowA`:V$rowA")
    $rngB = $ws.Range("F$rowB`:V$rowB")
    $valsA = $rngA.Value()
    $valsB = $rngB.Value()
    $rngA.Value = $valsB
    $rngB.Value = $valsA
}

function Rotate-RowData {
    # new(rowA) = old(rowB); new(rowB) = old(rowC); new(rowC) = old(rowA)
    param($rowA, $rowB, $rowC)
    $rngA = $ws.Range("F$rowA`:V$rowA")
    $rngB = $ws.Range("F$rowB`:V$rowB")
    $rngC = $ws.Range("F$rowC`:V$rowC")
    $valsA = $rngA.Value()
    $valsB = $rngB.Value()
    $valsC = $rngC.Value()
    $rngA.Value = $valsB
    $rngB.Value = $valsC
    $rngC.Value = $valsA
}

Swap-RowData 76 77
Rotate-RowData 84 85 86
Swap-RowData 99 100
Swap-RowData 102 103
Swap-RowData 105 107
Rotate-RowData 109 110 111
Swap-RowData 123 124
Rotate-RowData 142 143 144

# Append the new match row 147, copying row 146's formatting first.
$ws.Range("A146:V146").Copy() | Out-Null
$ws.Range("A147:V147").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Cells.Item(147, 1).Value = 146
$ws.Cells.Item(147, 2).Value = "italy"
$ws.Cells.Item(147, 3).Value = "serie-c-group-b"
$ws.Cells.Item(147, 4).Value = "2023-2024"
$ws.Cells.Item(147, 5).Value = 45257.86458333334
$ws.Cells.Item(147, 6).Value = "Recanatese"
$ws.Cells.Item(147, 7).Value = 1
$ws.Cells.Item(147, 8).Value = "Pineto"
$ws.Cells.Item(147, 9).Value = 1
$ws.Cells.Item(147, 10).Value = 2.12
$ws.Cells.Item(147, 11).Value = "23/11/2023 09:13"
$ws.Cells.Item(147, 12).Value = 2.13
$ws.Cells.Item(147, 13).Value = "27/11/2023 19:57"
$ws.Cells.Item(147, 14).Value = 2.99
$ws.Cells.Item(147, 15).Value = "23/11/2023 09:13"
$ws.Cells.Item(147, 16).Value = 3.1
$ws.Cells.Item(147, 17).Value = "27/11/2023 19:57"
$ws.Cells.Item(147, 18).Value = 3.38
$ws.Cells.Item(147, 19).Value = "23/11/2023 09:13"
$ws.Cells.Item(147, 20).Value = 3.78
$ws.Cells.Item(147, 21).Value = "27/11/2023 19:57"
$ws.Cells.Item(147, 22).Value = "https://www.betexplorer.com/football/italy/serie-c-group-b/recanatese-pineto/0vu3XrwE/"

Write-Output "Edit complete"
